$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.9155336666666667
$ws.Range("H2").Value = 2.746601
$ws.Range("I2").Value = 0.1890240037548773
$ws.Range("J2").Value = 0.1951995261655112
$ws.Range("M2").Value = 145.7087706666667
$ws.Range("N2").Value = 437.126312
$ws.Range("O2").Value = 0.5445232453600627
$ws.Range("P2").Value = 0.5461141113270247
$ws.Range("Q2").Value = 133.4012850739458
$ws.Range("R2").Value = 1200.611565665512
$ws.Range("S2").Value = 0.1029279639755585
$ws.Range("T2").Value = 0.1066012157633344
$ws.Range("G3").Value = 0.9155336666666667
$ws.Range("H3").Value = 2.746601
$ws.Range("I3").Value = 0.1890240037548773
$ws.Range("J3").Value = 0.1951995261655112
$ws.Range("O3").Value = 0.003349722554576428
$ws.Range("P3").Value = 0.003359509023117945
$ws.Range("Q3").Value = 0.8206395176503334
$ws.Range("R3").Value = 7.385755658853
$ws.Range("S3").Value = 0.0006331779687340519
$ws.Range("T3").Value = 0.0006557745694613822
$ws.Range("G4").Value = 0.9155336666666667
$ws.Range("H4").Value = 2.746601
$ws.Range("I4").Value = 0.1890240037548773
$ws.Range("J4").Value = 0.1951995261655112
$ws.Range("M4").Value = 51.59199533333333
$ws.Range("N4").Value = 154.775986
$ws.Range("O4").Value = 0.1928026748491032
$ws.Range("P4").Value = 0.1933659624890163
$ws.Range("Q4").Value = 47.23420865817622
$ws.Range("R4").Value = 425.107877923586
$ws.Range("S4").Value = 0.03644433353462726
$ws.Range("T4").Value = 0.03774494425439399
$ws.Range("G5").Value = 0.9155336666666667
$ws.Range("H5").Value = 2.746601
$ws.Range("I5").Value = 0.1890240037548773
$ws.Range("J5").Value = 0.1951995261655112
$ws.Range("M5").Value = 2.338518
$ws.Range("N5").Value = 4.677036
$ws.Range("O5").Value = 0.0087391953474509
$ws.Range("P5").Value = 0.005843151713055659
$ws.Range("Q5").Value = 2.140991959106
$ws.Range("R5").Value = 12.845951754636
$ws.Range("S5").Value = 0.001651917694171165
$ws.Range("T5").Value = 0.00114058044570166
$ws.Range("G6").Value = 0.9155336666666667
$ws.Range("H6").Value = 2.746601
$ws.Range("I6").Value = 0.1890240037548773
$ws.Range("J6").Value = 0.1951995261655112
$ws.Range("M6").Value = 67.05398933333333
$ws.Range("N6").Value = 201.161968
$ws.Range("O6").Value = 0.2505851618888069
$ws.Range("P6").Value = 0.2513172654477853
$ws.Range("Q6").Value = 61.39018471897423
$ws.Range("R6").Value = 552.511662470768
$ws.Range("S6").Value = 0.04736661058178636
$ws.Range("T6").Value = 0.04905701113261968
$ws.Range("I7").Value = 0.7156204889943075
$ws.Range("J7").Value = 0.7390002200311341
$ws.Range("M7").Value = 145.7087706666667
$ws.Range("N7").Value = 437.126312
$ws.Range("O7").Value = 0.5445232453600627
$ws.Range("P7").Value = 0.5461141113270247
$ws.Range("Q7").Value = 505.0400529071582
$ws.Range("R7").Value = 4545.360476164424
$ws.Range("S7").Value = 0.3896719911133353
$ws.Range("T7").Value = 0.4035784484327785
$ws.Range("I8").Value = 0.7156204889943075
$ws.Range("J8").Value = 0.7390002200311341
$ws.Range("O8").Value = 0.003349722554576428
$ws.Range("P8").Value = 0.003359509023117945
$ws.Range("S8").Value = 0.002397130092501245
$ws.Range("T8").Value = 0.002482677907280741
$ws.Range("I9").Value = 0.7156204889943075
$ws.Range("J9").Value = 0.7390002200311341
$ws.Range("M9").Value = 51.59199533333333
$ws.Range("N9").Value = 154.775986
$ws.Range("O9").Value = 0.1928026748491032
$ws.Range("P9").Value = 0.1933659624890163
$ws.Range("Q9").Value = 178.8226194862358
$ws.Range("R9").Value = 1609.403575376122
$ws.Range("S9").Value = 0.1379735444549257
$ws.Range("T9").Value = 0.142897488825915
$ws.Range("I10").Value = 0.7156204889943075
$ws.Range("J10").Value = 0.7390002200311341
$ws.Range("M10").Value = 2.338518
$ws.Range("N10").Value = 4.677036
$ws.Range("O10").Value = 0.0087391953474509
$ws.Range("P10").Value = 0.005843151713055659
$ws.Range("Q10").Value = 8.105519311162
$ws.Range("R10").Value = 48.63311586697201
$ws.Range("S10").Value = 0.00625394724795959
$ws.Range("T10").Value = 0.00431809040162343
$ws.Range("I11").Value = 0.7156204889943075
$ws.Range("J11").Value = 0.7390002200311341
$ws.Range("M11").Value = 67.05398933333333
$ws.Range("N11").Value = 201.161968
$ws.Range("O11").Value = 0.2505851618888069
$ws.Range("P11").Value = 0.2513172654477853
$ws.Range("Q11").Value = 232.4153183476818
$ws.Range("R11").Value = 2091.737865129136
$ws.Range("S11").Value = 0.1793238760855857
$ws.Range("T11").Value = 0.1857235144635362
$ws.Range("F12").Value = 0.6666666666666666
$ws.Range("G12").Value = 0.002153333333333333
$ws.Range("H12").Value = 0.00646
$ws.Range("I12").Value = 0.0004445840747369229
$ws.Range("J12").Value = 0.0004591088909634862
$ws.Range("M12").Value = 145.7087706666667
$ws.Range("N12").Value = 437.126312
$ws.Range("O12").Value = 0.5445232453600627
$ws.Range("P12").Value = 0.5461141113270247
$ws.Range("Q12").Value = 0.3137595528355556
$ws.Range("R12").Value = 2.82383597552
$ws.Range("S12").Value = 0.0002420863632111499
$ws.Range("T12").Value = 0.0002507258439908601
$ws.Range("F13").Value = 0.6666666666666666
$ws.Range("G13").Value = 0.002153333333333333
$ws.Range("H13").Value = 0.00646
$ws.Range("I13").Value = 0.0004445840747369229
$ws.Range("J13").Value = 0.0004591088909634862
$ws.Range("O13").Value = 0.003349722554576428
$ws.Range("P13").Value = 0.003359509023117945
$ws.Range("Q13").Value = 0.001930142486666667
$ws.Range("R13").Value = 0.01737128238
$ws.Range("S13").Value = 0.000001489233302551763
$ws.Range("T13").Value = 0.000001542380461785504
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.002153333333333333
$ws.Range("H14").Value = 0.00646
$ws.Range("I14").Value = 0.0004445840747369229
$ws.Range("J14").Value = 0.0004591088909634862
$ws.Range("M14").Value = 51.59199533333333
$ws.Range("N14").Value = 154.775986
$ws.Range("O14").Value = 0.1928026748491032
$ws.Range("P14").Value = 0.1933659624890163
$ws.Range("Q14").Value = 0.1110947632844444
$ws.Range("R14").Value = 0.99985286956
$ws.Range("S14").Value = 0.00008571699880459234
$ws.Range("T14").Value = 0.00008877603258841934
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.002153333333333333
$ws.Range("H15").Value = 0.00646
$ws.Range("I15").Value = 0.0004445840747369229
$ws.Range("J15").Value = 0.0004591088909634862
$ws.Range("M15").Value = 2.338518
$ws.Range("N15").Value = 4.677036
$ws.Range("O15").Value = 0.0087391953474509
$ws.Range("P15").Value = 0.005843151713055659
$ws.Range("Q15").Value = 0.005035608760000001
$ws.Range("R15").Value = 0.03021365256
$ws.Range("S15").Value = 0.00000388530707749168
$ws.Range("T15").Value = 0.000002682642902712378
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.002153333333333333
$ws.Range("H16").Value = 0.00646
$ws.Range("I16").Value = 0.0004445840747369229
$ws.Range("J16").Value = 0.0004591088909634862
$ws.Range("M16").Value = 67.05398933333333
$ws.Range("N16").Value = 201.161968
$ws.Range("O16").Value = 0.2505851618888069
$ws.Range("P16").Value = 0.2513172654477853
$ws.Range("Q16").Value = 0.1443895903644445
$ws.Range("R16").Value = 1.29950631328
$ws.Range("S16").Value = 0.0001114061723411372
$ws.Range("T16").Value = 0.0001153819910197088
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 0.459699
$ws.Range("H17").Value = 0.9193979999999999
$ws.Range("I17").Value = 0.09491092317607834
$ws.Range("J17").Value = 0.06534114491239122
$ws.Range("M17").Value = 145.7087706666667
$ws.Range("N17").Value = 437.126312
$ws.Range("O17").Value = 0.5445232453600627
$ws.Range("P17").Value = 0.5461141113270247
$ws.Range("Q17").Value = 66.98217616669599
$ws.Range("R17").Value = 401.893057000176
$ws.Range("S17").Value = 0.05168120390795777
$ws.Range("T17").Value = 0.03568372128692087
$ws.Range("E18").Value = 2
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 0.459699
$ws.Range("H18").Value = 0.9193979999999999
$ws.Range("I18").Value = 0.09491092317607834
$ws.Range("J18").Value = 0.06534114491239122
$ws.Range("O18").Value = 0.003349722554576428
$ws.Range("P18").Value = 0.003359509023117945
$ws.Range("Q18").Value = 0.412051658349
$ws.Range("R18").Value = 2.472309950094
$ws.Range("S18").Value = 0.0003179252600385803
$ws.Range("T18").Value = 0.0002195141659140355
$ws.Range("E19").Value = 2
$ws.Range("F19").Value = 1
$ws.Range("G19").Value = 0.459699
$ws.Range("H19").Value = 0.9193979999999999
$ws.Range("I19").Value = 0.09491092317607834
$ws.Range("J19").Value = 0.06534114491239122
$ws.Range("M19").Value = 51.59199533333333
$ws.Range("N19").Value = 154.775986
$ws.Range("O19").Value = 0.1928026748491032
$ws.Range("P19").Value = 0.1933659624890163
$ws.Range("Q19").Value = 23.716788662738
$ws.Range("R19").Value = 142.300731976428
$ws.Range("S19").Value = 0.01829907986074564
$ws.Range("T19").Value = 0.01263475337611882
$ws.Range("E20").Value = 2
$ws.Range("F20").Value = 1
$ws.Range("G20").Value = 0.459699
$ws.Range("H20").Value = 0.9193979999999999
$ws.Range("I20").Value = 0.09491092317607834
$ws.Range("J20").Value = 0.06534114491239122
$ws.Range("M20").Value = 2.338518
$ws.Range("N20").Value = 4.677036
$ws.Range("O20").Value = 0.0087391953474509
$ws.Range("P20").Value = 0.005843151713055659
$ws.Range("Q20").Value = 1.075014386082
$ws.Range("R20").Value = 4.300057544328
$ws.Range("S20").Value = 0.0008294450982426535
$ws.Range("T20").Value = 0.0003817982228278568
$ws.Range("E21").Value = 2
$ws.Range("F21").Value = 1
$ws.Range("G21").Value = 0.459699
$ws.Range("H21").Value = 0.9193979999999999
$ws.Range("I21").Value = 0.09491092317607834
$ws.Range("J21").Value = 0.06534114491239122
$ws.Range("M21").Value = 67.05398933333333
$ws.Range("N21").Value = 201.161968
$ws.Range("O21").Value = 0.2505851618888069
$ws.Range("P21").Value = 0.2513172654477853
$ws.Range("Q21").Value = 30.824651842544
$ws.Range("R21").Value = 184.947911055264
$ws.Range("S21").Value = 0.0237832690490937
$ws.Range("T21").Value = 0.01642135786060963
